$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.750.07'
$ws.Range("E2").Value = '  +6.37%  '
$ws.Range("D3").Value = '1.738.13'
$ws.Range("E3").Value = '  +5.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.46'
$ws.Range("E5").Value = '  +4.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5459'
$ws.Range("E6").Value = '  +3.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2764'
$ws.Range("E8").Value = '  +3.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06724'
$ws.Range("E9").Value = '  +5.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.94'
$ws.Range("E10").Value = '  +6.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07782'
$ws.Range("E11").Value = '  +1.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.682'
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("D13").Value = '1.750.13'
$ws.Range("E13").Value = '  +3.67%  '
$ws.Range("D14").Value = '1.976.86'
$ws.Range("E14").Value = '  +5.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5969'
$ws.Range("E15").Value = '  +6.30%  '
$ws.Range("D16").Value = '0.0₅8393'
$ws.Range("E16").Value = '  +2.07%  '
$ws.Range("E17").Value = '  +5.62%  '
$ws.Range("D18").Value = '27.750.44'
$ws.Range("E18").Value = '  +6.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '224.77'
$ws.Range("E19").Value = '  +17.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.827'
$ws.Range("E20").Value = '  +3.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.88'
$ws.Range("E22").Value = '  +5.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.225'
$ws.Range("E23").Value = '  +4.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.005'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.94'
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.714'
$ws.Range("E26").Value = '  +14.42%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1249'
$ws.Range("E27").Value = '  +4.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.450'
$ws.Range("E28").Value = '  +2.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '17.14'
$ws.Range("E29").Value = '  +7.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05669'
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("E31").Value = '  +3.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.687'
$ws.Range("E32").Value = '  +5.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.513'
$ws.Range("E33").Value = '  +4.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.679'
$ws.Range("E34").Value = '  +6.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9768'
$ws.Range("E35").Value = '  +3.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.858'
$ws.Range("E36").Value = '  +2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.452'
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5950'
$ws.Range("E38").Value = '  +3.14%  '
$ws.Range("E39").Value = '  +4.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.885'
$ws.Range("E40").Value = '  -1.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8493'
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("D42").Value = '1.047.50'
$ws.Range("E42").Value = '  +2.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.004'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.80'
$ws.Range("E44").Value = '  +0.36%  '
$ws.Range("D45").Value = '1.882.87'
$ws.Range("E45").Value = '  +5.13%  '
$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").Value = '  +8.30%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.30'
$ws.Range("E47").Value = '  +1.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.283'
$ws.Range("E48").Value = '  +3.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4436'
$ws.Range("E49").Value = '  +2.14%  '
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05318'
$ws.Range("E51").Value = '  -0.63%  '
